# Correct the "PACS Submit status" value in O2 from "Successfully verified"
# to "Successfully Verified" (capital V).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("O2").Value = "Successfully Verified"
